# Add a new worksheet "Time Consup. Anlys.S." after the existing
# "Incomplete US labelling" sheet, and populate it with the base timing
# data that Alex and Lukas have already analysed.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "Time Consup. Anlys.S."

# --- header row ---------------------------------------------------------
$headers = New-Object 'object[,]' 1,8
$headers[0,0] = "Dataset"
$headers[0,1] = "Run Count"
$headers[0,2] = "Model Version"
$headers[0,3] = "Threading Enabled"
$headers[0,4] = "Nanoseconds"
$headers[0,5] = "Milliseconds"
$headers[0,6] = "Seconds"
$headers[0,7] = "Minutes"
$ws2.Range("A1:H1").Value = $headers

# --- data rows (8 columns, incl. "true" placeholder for col D) ---------
$data = @(
    @("g03",0,"gpt-3.5-turbo","true",38512547600,38512.5476,38.5125476,0.6418757933333333),
    @("g04",0,"gpt-3.5-turbo","true",40979311400,40979.3114,40.9793114,0.6829885233333334),
    @("g08",0,"gpt-3.5-turbo","true",34931253200,34931.2532,34.9312532,0.5821875533333334),
    @("g10",0,"gpt-3.5-turbo","true",36928688300,36928.6883,36.9286883,0.6154781383333333),
    @("g11",0,"gpt-3.5-turbo","true",62972494100,62972.4941,62.9724941,1.049541568333333),
    @("g14",0,"gpt-3.5-turbo","true",28285799500,28285.7995,28.2857995,0.4714299916666667),
    @("g16",0,"gpt-3.5-turbo","true",62265483300,62265.4833,62.2654833,1.037758055),
    @("g18",0,"gpt-3.5-turbo","true",69495710600,69495.7106,69.4957106,1.158261843333333),
    @("g19",0,"gpt-3.5-turbo","true",60467353300,60467.3533,60.4673533,1.007789221666667),
    @("g21",0,"gpt-3.5-turbo","true",35529855900,35529.8559,35.5298559,0.592164265),
    @("g22",0,"gpt-3.5-turbo","true",43263083000,43263.083,43.263083,0.7210513833333334),
    @("g23",0,"gpt-3.5-turbo","true",40747012900,40747.0129,40.7470129,0.6791168816666667),
    @("g24",0,"gpt-3.5-turbo","true",32971244800,32971.2448,32.9712448,0.5495207466666667),
    @("g25",0,"gpt-3.5-turbo","true",78224568600,78224.5686,78.2245686,1.30374281),
    @("g26",0,"gpt-3.5-turbo","true",17749166200,17749.1662,17.7491662,0.2958194366666667),
    @("g27",0,"gpt-3.5-turbo","true",20200447000,20200.447,20.200447,0.3366741166666667),
    @("g28",0,"gpt-3.5-turbo","true",32354950100,32354.9501,32.3549501,0.5392491683333334)
)

$rowCount = $data.Length
$colCount = 8
$block = New-Object 'object[,]' $rowCount,$colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $block[$r,$c] = $data[$r][$c]
    }
}
$ws2.Range("A2:H18").Value = $block

# "Threading Enabled" values must stay literal text ("true"), not get
# auto-converted into a boolean by the Value setter. Write them with a
# leading quote-prefix (keeps them as text) then strip the quote-prefix
# number format back off by pasting the plain-text format from a normal
# data cell in the same row band, so the cells end up styled exactly
# like the rest of the data rows.
$quoted = New-Object 'object[,]' 17,1
for ($i = 0; $i -lt 17; $i++) { $quoted[$i,0] = "'true" }
$ws2.Range("D2:D18").Value = $quoted

$ws2.Range("C2").Copy()
$ws2.Range("D2:D18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- header / data cell styling -----------------------------------------
# Match the header/body styles already used on "Incomplete US labelling"
# (bold centered header row, wrapped-text body rows).
$ws1.Range("A1").Copy()
$ws2.Range("A1:H1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws2.Range("A2:H18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- column widths --------------------------------------------------------
$widths = @(12.666666666666666,15.666666666666666,21.666666666666668,27.666666666666668,18.666666666666668,20.166666666666668,12.666666666666666,12.666666666666666)
for ($c = 1; $c -le 8; $c++) {
    $ws2.Columns.Item($c).ColumnWidth = $widths[$c-1]
}

# --- freeze header row, matching the other sheet -------------------------
$ws2.Activate()
$ws2.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws2.Range("A1").Select()

# --- outline properties (summary rows below / summary cols to the right) -
$ws2.Outline.SummaryRow = 1
$ws2.Outline.SummaryColumn = -4152

# --- autofilter + matching defined name -----------------------------------
$ws2.Range("A1:H1").AutoFilter()
$fdbName = $ws2.Names.Add("_xlnm._FilterDatabase", "='Time Consup. Anlys.S.'!`$A`$1:`$H`$1")
$fdbName.Visible = $false

# Leave the workbook with the original sheet active, as it was before.
$ws1.Activate()
